$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the "The ELSEIF" / "2" / " " runs in paragraph 2 into a
#    single run "The ELSEIF2 ", while leaving the following "paragraph"
#    and "." runs untouched (still two separate runs).
#
#    The COM shim merges every run from the edited position through to
#    the end of the paragraph whenever text is changed, so we first
#    split the paragraph right after the " " run, rebuild the merged
#    run in isolation, and then rejoin the two paragraphs by deleting
#    the paragraph mark we introduced. This keeps "paragraph" and "."
#    as separate runs exactly like before.
# ---------------------------------------------------------------------

$splitPoint = $d.Content.Find.Execute("The ELSEIF2 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeRange = $d.Content
$mergeRange.Find.Execute("The ELSEIF2 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeStart = $mergeRange.Start
$mergeEnd = $mergeRange.End

$splitRange = $d.Range($mergeEnd, $mergeEnd)
$splitRange.InsertParagraphAfter()

$oldRunsRange = $d.Range($mergeStart, $mergeEnd)
$oldRunsRange.Delete()

$newRunRange = $d.Range($mergeStart, $mergeStart)
$newRunRange.InsertAfter("The ELSEIF2 ")

$pmark = $d.Range($mergeEnd, $mergeEnd + 1)
$pmark.Delete()

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the last (empty) paragraph to the
#    very start of the document/first paragraph.
#
#    Adding a zero-length bookmark directly at position 0 is special
#    cased by the engine (it balloons to cover the whole first
#    paragraph), so we insert a throwaway character at position 0,
#    anchor the bookmark to the position right after it (which behaves
#    correctly), and then delete the throwaway character again - the
#    zero-length bookmark collapses back down to position 0.
# ---------------------------------------------------------------------

$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")

$bookmarkRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$d.Range(0, 1).Delete()
